$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on ambiguous numeric-looking price cells in column D
# so Excel does not auto-convert them to numeric values (preserves original
# inline-string formatting of the source data, e.g. '8.096', '1.002').
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values from the crypto price refresh
$ws.Range("D2").Value = '27.379.04'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '1.817.36'
$ws.Range("E3").Value = '  -3.11%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.87%  '
$ws.Range("D5").Value = '331.73'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '0.4566'
$ws.Range("E7").Value = '  -2.05%  '
$ws.Range("D8").Value = '0.3807'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").Value = '45.96'
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").Value = '0.07849'
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").Value = '0.9577'
$ws.Range("E11").Value = '  -5.08%  '
$ws.Range("D12").Value = '20.96'
$ws.Range("E12").Value = '  -4.16%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.824.86'
$ws.Range("E13").Value = '  -3.37%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.839'
$ws.Range("E14").Value = '  -2.40%  '
$ws.Range("D15").Value = '7.055'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '89.34'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '0.06587'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '0.00001019'
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("D20").Value = '17.10'
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").Value = '27.364.49'
$ws.Range("D23").Value = '5.281'
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("D24").Value = '10.80'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '2.258'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '2.032.70'
$ws.Range("E26").Value = '  -3.53%  '
$ws.Range("D27").Value = '155.56'
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = '19.28'
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = '2.042'
$ws.Range("E29").Value = '  -5.14%  '
$ws.Range("D30").Value = '5.241'
$ws.Range("E30").Value = '  -4.43%  '
$ws.Range("D31").Value = '117.55'
$ws.Range("E31").Value = '  -3.41%  '
$ws.Range("D32").Value = '0.09286'
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").Value = '0.9311'
$ws.Range("E33").Value = '  -4.96%  '
$ws.Range("D34").Value = '3.566'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = '5.211'
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("D36").Value = '1.313'
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").Value = '0.05904'
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("D38").Value = '0.02177'
$ws.Range("E38").Value = '  -2.84%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '8.096'
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").Value = '1.139'
$ws.Range("E41").Value = '  -5.17%  '
$ws.Range("D42").Value = '0.5729'
$ws.Range("E42").Value = '  -3.92%  '
$ws.Range("D43").Value = '0.1812'
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("D44").Value = '9.914'
$ws.Range("E44").Value = '  -4.54%  '
$ws.Range("D45").Value = '1.267'
$ws.Range("E45").Value = '  +1.35%  '
$ws.Range("D46").Value = '11.84'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").Value = '0.5383'
$ws.Range("E47").Value = '  -4.67%  '
$ws.Range("D48").Value = '1.864'
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("D49").Value = '0.06572'
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("D50").Value = '109.57'
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  -33.67%  '
